$wb = $excel.ActiveWorkbook

# Add the new "FlightFinder" sheet after the existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "FlightFinder"

# Populate FlightFinder data (row 1)
$ws2.Range("A1").Value = 2
$ws2.Range("B1").Value = "London"
$ws2.Range("C1").Value = "December"
$ws2.Range("D1").Value = 10
$ws2.Range("E1").Value = "Paris"
$ws2.Range("F1").Value = "December"
$ws2.Range("G1").Value = 30
$ws2.Range("H1").Value = "Unified Airlines"

# Column H width (bestFit/customWidth, resolves to stored width 15)
$ws2.Columns.Item(8).ColumnWidth = 14.14

# Select H1 on FlightFinder and make it the active sheet (tabSelected)
$ws2.Activate()
$ws2.Range("H1").Select()
